# Update "countries & provincias Spain" data dump
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 15:22"

# --- Straightforward per-country stat refreshes (no reordering) ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 1064832
$ws.Range("C4").Value = 638
$ws.Range("E4").Value = 855679

# España (row 5)
$ws.Range("B5").Value = 239639
$ws.Range("C5").Value = 2740
$ws.Range("D5").Value = 137984
$ws.Range("E5").Value = 77112
$ws.Range("F5").Value = 2676
$ws.Range("G5").Value = 268
$ws.Range("H5").Value = 24543

# Suiza (row 20)
$ws.Range("E20").Value = 4749
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 1737

# Arabia Saudita (row 23)
$ws.Range("B23").Value = 22753
$ws.Range("C23").Value = 1351
$ws.Range("D23").Value = 3163
$ws.Range("E23").Value = 19428
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 162

# Pakistan (row 28)
$ws.Range("B28").Value = 16117
$ws.Range("C28").Value = 592
$ws.Range("D28").Value = 4105
$ws.Range("E28").Value = 11654

# Serbia (row 42)
$ws.Range("B42").Value = 9009
$ws.Range("C42").Value = 285
$ws.Range("D42").Value = 1343
$ws.Range("E42").Value = 7487
$ws.Range("G42").Value = 6
$ws.Range("H42").Value = 179

# Finlandia (row 54)
$ws.Range("D54").Value = 3000
$ws.Range("E54").Value = 1784

# Croacia (row 68)
$ws.Range("F68").Value = 20

# Uzbekistan (row 71)
$ws.Range("D71").Value = 1103
$ws.Range("E71").Value = 905

# Sri Lanka (row 103)
$ws.Range("B103").Value = 653
$ws.Range("C103").Value = 4
$ws.Range("E103").Value = 507

# San Marino (row 108)
$ws.Range("B108").Value = 569
$ws.Range("C108").Value = 6
$ws.Range("D108").Value = 78
$ws.Range("E108").Value = 450

# --- Reorder: Kenia now overtakes El Salvador (rows 119-120) ---
$ws.Range("A119").Value = "Kenia"
$ws.Range("B119").Value = 396
$ws.Range("C119").Value = 12
$ws.Range("D119").Value = 144
$ws.Range("E119").Value = 235
$ws.Range("F119").Value = 2
$ws.Range("G119").Value = 2
$ws.Range("H119").Value = 17

$ws.Range("A120").Value = "El Salvador"
$ws.Range("B120").Value = 395
$ws.Range("C120").Value = 18
$ws.Range("D120").Value = 118
$ws.Range("E120").Value = 268
$ws.Range("F120").Value = 3
$ws.Range("G120").Value = 0
$ws.Range("H120").Value = 9

# --- Reorder: Togo now overtakes Trinidad yTobago, Cabo Verde, Bermudas (rows 147-150) ---
$ws.Range("A147").Value = "Togo"
$ws.Range("B147").Value = 116
$ws.Range("C147").Value = 7
$ws.Range("D147").Value = 65
$ws.Range("E147").Value = 42
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 2
$ws.Range("H147").Value = 9

$ws.Range("A148").Value = "Trinidad yTobago"
$ws.Range("B148").Value = 116
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 72
$ws.Range("E148").Value = 36
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 8

$ws.Range("A149").Value = "Cabo Verde"
$ws.Range("B149").Value = 113
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 2
$ws.Range("E149").Value = 110
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 1

$ws.Range("A150").Value = "Bermudas"
$ws.Range("B150").Value = 111
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 48
$ws.Range("E150").Value = 57
$ws.Range("F150").Value = 10
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 6
